$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting (bold, centered, bordered) from H1 into new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row additions
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows for new columns I (I0) and J (IF)
$data = @(
    @(9, 9),
    @(8, 9),
    @(7, 7),
    @(7, 7),
    @(6, 8),
    @(7, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
